# Insert two new weekly-report rows right before the current row 218.
# This pushes the existing rows 218-330 down to 220-332 (Excel carries every
# column along with the row, so no other edits are needed for those rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(218).Insert()
$ws.Rows.Item(218).Insert()

# The newly inserted rows 218/219 are blank except for the date-cell style
# that Excel carried over from the insert point. Clone the (now shifted)
# row 220 into both of them so every constant column (Mercado, Region,
# Categoria, etc.) is populated exactly like the rest of the sheet.
$ws.Range("A220:R220").Copy()
$ws.Range("A218:R218").PasteSpecial()
$ws.Range("A220:R220").Copy()
$ws.Range("A219:R219").PasteSpecial()

# Row 218: new "Segunda" quality record for 2022-03-16 (serial 44636).
$ws.Range("D218").Value = 44636
$ws.Range("I218").Value = "Segunda"
$ws.Range("J218").Value = 1200
$ws.Range("K218").Value = 400
$ws.Range("L218").Value = 500
$ws.Range("M218").Value = 450
$ws.Range("P218").Value = 450

# Row 219: new "Tercera" quality record for 2022-03-16 (serial 44636).
$ws.Range("D219").Value = 44636
$ws.Range("I219").Value = "Tercera"
$ws.Range("J219").Value = 1200
$ws.Range("K219").Value = 300
$ws.Range("L219").Value = 350
$ws.Range("M219").Value = 325
$ws.Range("P219").Value = 325
